# Updates the "想去人数" (F column) figures across all four sheets to
# match the refreshed data snapshot (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 36
$ws1.Range("F3").Value = 1312
$ws1.Range("F4").Value = 13058
$ws1.Range("F5").Value = 745
$ws1.Range("F8").Value = 62
$ws1.Range("F10").Value = 1893
$ws1.Range("F11").Value = 42
$ws1.Range("F13").Value = 8580
$ws1.Range("F15").Value = 217
$ws1.Range("F17").Value = 359
$ws1.Range("F18").Value = 228
$ws1.Range("F19").Value = 307
$ws1.Range("F20").Value = 146
$ws1.Range("F22").Value = 31
$ws1.Range("F23").Value = 228
$ws1.Range("F24").Value = 268
$ws1.Range("F25").Value = 1324
$ws1.Range("F27").Value = 71
$ws1.Range("F28").Value = 104

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 4456
$ws2.Range("F6").Value = 172
$ws2.Range("F8").Value = 21
$ws2.Range("F9").Value = 77
$ws2.Range("F10").Value = 77
$ws2.Range("F11").Value = 372
$ws2.Range("F16").Value = 14
$ws2.Range("F17").Value = 16
$ws2.Range("F19").Value = 14

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 892
$ws3.Range("F3").Value = 4324
$ws3.Range("F4").Value = 9

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 892
$ws4.Range("F3").Value = 36
$ws4.Range("F6").Value = 1312
$ws4.Range("F7").Value = 13058
$ws4.Range("F9").Value = 745
$ws4.Range("F10").Value = 4324
$ws4.Range("F13").Value = 62
$ws4.Range("F15").Value = 1893
$ws4.Range("F16").Value = 42
$ws4.Range("F18").Value = 9
$ws4.Range("F19").Value = 8584
$ws4.Range("F21").Value = 4456
$ws4.Range("F22").Value = 217
$ws4.Range("F23").Value = 172
$ws4.Range("F24").Value = 172
$ws4.Range("F27").Value = 21
$ws4.Range("F28").Value = 77
$ws4.Range("F29").Value = 77
$ws4.Range("F30").Value = 372
$ws4.Range("F31").Value = 359
$ws4.Range("F33").Value = 228
$ws4.Range("F34").Value = 307
$ws4.Range("F35").Value = 146
$ws4.Range("F37").Value = 31
$ws4.Range("F39").Value = 228
$ws4.Range("F42").Value = 268
$ws4.Range("F43").Value = 1324
$ws4.Range("F44").Value = 14
$ws4.Range("F46").Value = 71
$ws4.Range("F47").Value = 104
$ws4.Range("F48").Value = 16
$ws4.Range("F50").Value = 14
